$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B values for rows 4-73 per the naive forecaster bugfix
$ws.Range("B4").Value = 0.9623985945846414
$ws.Range("B5").Value = -0.2625025240627394
$ws.Range("B6").Value = -2.717991563576632
$ws.Range("B7").Value = -0.8018087606025261
$ws.Range("B8").Value = -2.540735588184205
$ws.Range("B9").Value = 0.4798059474883161
$ws.Range("B10").Value = 2.996743291460049
$ws.Range("B11").Value = 1.766782973262963
$ws.Range("B12").Value = 1.792707332192279
$ws.Range("B13").Value = 1.599999999999994
$ws.Range("B14").Value = 1.983944815439088
$ws.Range("B15").Value = 1.027205282249909
$ws.Range("B16").Value = 1.892417316869597
$ws.Range("B17").Value = 1.727088565964991
$ws.Range("B18").Value = 0.4000000000000057
$ws.Range("B19").Value = 0.5999999999999943
$ws.Range("B20").Value = 1.295244683175738
$ws.Range("B21").Value = 0.5000000000000142
$ws.Range("B22").Value = -0.1007444012410019
$ws.Range("B23").Value = 0.399975276622385
$ws.Range("B24").Value = 2
$ws.Range("B25").Value = 1.1
$ws.Range("B26").Value = 1.292635181922734
$ws.Range("B27").Value = 1.682020243440505
$ws.Range("B28").Value = 1.012497979540356
$ws.Range("B29").Value = 1.634557118349079
$ws.Range("B30").Value = 0.8971548841028039
$ws.Range("B31").Value = 1.102490924221428
$ws.Range("B32").Value = 0.9984536597660991
$ws.Range("B33").Value = 1.495216092286043
$ws.Range("B34").Value = 1.345302727311861
$ws.Range("B35").Value = 1.428495556385869
$ws.Range("B36").Value = 0.7394362573787987
$ws.Range("B37").Value = 1.299999999999983
$ws.Range("B38").Value = 1.233394657233262
$ws.Range("B39").Value = 1.297884859680252
$ws.Range("B40").Value = 1.637264503497377
$ws.Range("B41").Value = 1.331979115623398
$ws.Range("B42").Value = 1.795921598270084
$ws.Range("B43").Value = 1.327129713066284
$ws.Range("B44").Value = 1.15883627336575
$ws.Range("B45").Value = 1.076058203620576
$ws.Range("B46").Value = 0.3
$ws.Range("B47").Value = 2.40590956953757
$ws.Range("B48").Value = -0.4
$ws.Range("B49").Value = 0.8963263664365542
$ws.Range("B50").Value = 1.055400706275506
$ws.Range("B51").Value = -0.4303992348575321
$ws.Range("B52").Value = -14.5
$ws.Range("B53").Value = 7.96955251685678
$ws.Range("B54").Value = 2.117788110998191
$ws.Range("B55").Value = -0.9861240056009706
$ws.Range("B56").Value = 2.876944405321424
$ws.Range("B57").Value = 0.4364757668776207
$ws.Range("B58").Value = 0.6542354095451515
$ws.Range("B59").Value = -1.358640149334988
$ws.Range("B60").Value = 1.767346889326234
$ws.Range("B61").Value = 0.1973819540654631
$ws.Range("B62").Value = -2.016521230865749
$ws.Range("B63").Value = -1.429923541452922
$ws.Range("B64").Value = 0.1801827438520291
$ws.Range("B65").Value = 0.6692453970872521
$ws.Range("B66").Value = -1.33243152085096
$ws.Range("B67").Value = -0.09168777270478756
$ws.Range("B68").Value = 0.8481739611978583
$ws.Range("B69").Value = -0.2604183589432552
$ws.Range("B70").Value = 0.3093370292089048
$ws.Range("B71").Value = -0.03180050048325711
$ws.Range("B72").Value = 0.04735640278761366
$ws.Range("B73").Value = 0.2461857363876589

# Remove now-obsolete forecast rows 74-82 (dimension shrinks from B82 to B73)
$ws.Range("A74:B82").EntireRow.Delete()

